$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 74
$ws.Range("I4").Value = 77
$ws.Range("K4").Value = 77
$ws.Range("M4").Value = 37
$ws.Range("H18").Value = 822.6667
$ws.Range("I18").Value = 357.2
$ws.Range("J18").Value = 3150
$ws.Range("K18").Value = 357.2
$ws.Range("L18").Value = 3150
$ws.Range("M18").Value = -73.19999999999999
$ws.Range("N18").Value = -3718
$ws.Range("H103").Value = 1195.625
$ws.Range("I103").Value = 1148.6666
$ws.Range("J103").Value = 1336.5
$ws.Range("K103").Value = 3445.9998
$ws.Range("L103").Value = 4009.5
$ws.Range("M103").Value = -2859.9998
$ws.Range("N103").Value = -5181.5
$ws.Range("H106").Value = 1991
$ws.Range("I106").Value = 789.2
$ws.Range("K106").Value = 789.2
$ws.Range("M106").Value = -158.2
$ws.Range("H112").Value = 3420.4348
$ws.Range("I112").Value = 1729.7
$ws.Range("J112").Value = 4721
$ws.Range("K112").Value = 5189.1
$ws.Range("L112").Value = 14163
$ws.Range("M112").Value = -4081.1
$ws.Range("N112").Value = -16379
$ws.Range("H113").Value = 5495.4
$ws.Range("I113").Value = 5936.8887
$ws.Range("J113").Value = 4833.1665
$ws.Range("K113").Value = 5936.8887
$ws.Range("L113").Value = 4833.1665
$ws.Range("M113").Value = -2682.8887
$ws.Range("N113").Value = -11341.1665
$ws.Range("H132").Value = 5841.85
$ws.Range("I132").Value = 6106.875
$ws.Range("J132").Value = 4781.75
$ws.Range("K132").Value = 18320.625
$ws.Range("L132").Value = 14345.25
$ws.Range("M132").Value = -15790.625
$ws.Range("N132").Value = -19405.25
$ws.Range("H137").Value = 6066.385
$ws.Range("I137").Value = 1393.9231
$ws.Range("J137").Value = 10738.846
$ws.Range("K137").Value = 4181.7693
$ws.Range("L137").Value = 32216.538
$ws.Range("M137").Value = -1631.7693
$ws.Range("N137").Value = -37316.538
$ws.Range("H138").Value = 11504.66
$ws.Range("I138").Value = 3699.5
$ws.Range("J138").Value = 12991.357
$ws.Range("K138").Value = 11098.5
$ws.Range("L138").Value = 38974.071
$ws.Range("M138").Value = -5958.5
$ws.Range("N138").Value = -49254.071

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 33372686
$ws.Range("I74").Value = 55619120
$ws.Range("J74").Value = 3033.1667
$ws.Range("K74").Value = 55619120
$ws.Range("L74").Value = 3033.1667
$ws.Range("M74").Value = -55618246
$ws.Range("N74").Value = -4781.1667
$ws.Range("H77").Value = 33372686
$ws.Range("I77").Value = 55619120
$ws.Range("J77").Value = 3033.1667
$ws.Range("K77").Value = 278095600
$ws.Range("L77").Value = 15165.8335
$ws.Range("M77").Value = -278091232
$ws.Range("N77").Value = -23901.8335
$ws.Range("H122").Value = 7094337.5
$ws.Range("I122").Value = 1726.2433
$ws.Range("K122").Value = 5178.7299
$ws.Range("M122").Value = -2728.7299
$ws.Range("H132").Value = 52711416
$ws.Range("I132").Value = 9491.299999999999
$ws.Range("J132").Value = 111269110
$ws.Range("K132").Value = 28473.9
$ws.Range("L132").Value = 333807330
$ws.Range("M132").Value = -25943.9
$ws.Range("N132").Value = -333812390

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 947.5
$ws.Range("I16").Value = 1028.3334
$ws.Range("J16").Value = 802
$ws.Range("K16").Value = 1028.3334
$ws.Range("L16").Value = 802
$ws.Range("M16").Value = -741.3334
$ws.Range("N16").Value = -1376
$ws.Range("H31").Value = 27780154
$ws.Range("I31").Value = 2635.5715
$ws.Range("J31").Value = 40324840
$ws.Range("K31").Value = 2635.5715
$ws.Range("L31").Value = 40324840
$ws.Range("M31").Value = -2340.5715
$ws.Range("N31").Value = -40325430
$ws.Range("H34").Value = 27780154
$ws.Range("I34").Value = 2635.5715
$ws.Range("J34").Value = 40324840
$ws.Range("K34").Value = 2635.5715
$ws.Range("L34").Value = 40324840
$ws.Range("M34").Value = -2433.5715
$ws.Range("N34").Value = -40325244
$ws.Range("H99").Value = 8428.571
$ws.Range("I99").Value = 7000
$ws.Range("J99").Value = 8666.666999999999
$ws.Range("K99").Value = 7000
$ws.Range("L99").Value = 8666.666999999999
$ws.Range("M99").Value = -5502
$ws.Range("N99").Value = -11662.667
$ws.Range("H113").Value = 947.5
$ws.Range("I113").Value = 1028.3334
$ws.Range("J113").Value = 802
$ws.Range("K113").Value = 1028.3334
$ws.Range("L113").Value = 802
$ws.Range("M113").Value = 1141.6666
$ws.Range("N113").Value = -5142
$ws.Range("H122").Value = 2393925
$ws.Range("I122").Value = 1243.5834
$ws.Range("J122").Value = 5265142.5
$ws.Range("K122").Value = 3730.7502
$ws.Range("L122").Value = 15795427.5
$ws.Range("M122").Value = -1280.7502
$ws.Range("N122").Value = -15800327.5
$ws.Range("H126").Value = 8428.571
$ws.Range("I126").Value = 7000
$ws.Range("J126").Value = 8666.666999999999
$ws.Range("K126").Value = 21000
$ws.Range("L126").Value = 26000.001
$ws.Range("M126").Value = -18530
$ws.Range("N126").Value = -30940.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 18140.934
$ws.Range("I56").Value = 18140.934
$ws.Range("K56").Value = 18140.934
$ws.Range("M56").Value = -17610.934
$ws.Range("H68").Value = 7234.9614
$ws.Range("I68").Value = 6872.5
$ws.Range("J68").Value = 7396.0557
$ws.Range("K68").Value = 20617.5
$ws.Range("L68").Value = 22188.1671
$ws.Range("M68").Value = -19806.5
$ws.Range("N68").Value = -23810.1671
$ws.Range("H71").Value = 7234.9614
$ws.Range("I71").Value = 6872.5
$ws.Range("J71").Value = 7396.0557
$ws.Range("K71").Value = 61852.5
$ws.Range("L71").Value = 66564.5013
$ws.Range("M71").Value = -57796.5
$ws.Range("N71").Value = -74676.5013
$ws.Range("H107").Value = 2063.3044
$ws.Range("J107").Value = 2222.95
$ws.Range("L107").Value = 6668.849999999999
$ws.Range("N107").Value = -10508.85
$ws.Range("H129").Value = 2445.913
$ws.Range("I129").Value = 927.1429000000001
$ws.Range("J129").Value = 3110.375
$ws.Range("K129").Value = 2781.4287
$ws.Range("L129").Value = 9331.125
$ws.Range("M129").Value = 2218.5713
$ws.Range("N129").Value = -19331.125
$ws.Range("H131").Value = 19613878
$ws.Range("I131").Value = 1209.625
$ws.Range("J131").Value = 37047360
$ws.Range("K131").Value = 3628.875
$ws.Range("L131").Value = 111142080
$ws.Range("M131").Value = 1411.125
$ws.Range("N131").Value = -111152160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2557.3062
$ws.Range("I102").Value = 1804.4193
$ws.Range("J102").Value = 3853.9443
$ws.Range("K102").Value = 1804.4193
$ws.Range("L102").Value = 3853.9443
$ws.Range("M102").Value = -182.4193
$ws.Range("N102").Value = -7097.9443
$ws.Range("H132").Value = 5727.7144
$ws.Range("I132").Value = 4998.1177
$ws.Range("K132").Value = 14994.3531
$ws.Range("M132").Value = -12464.3531
